# Separate columns for vacant and unavailable
#
# The totals row (row 6) had its formula-placeholder strings off-by-one
# relative to the columns they sit above (e.g. the cell in column E summed
# column D). This lines each total back up with its own column, frees up
# column D's SUMIF to key off column D itself (instead of C), and switches
# the duplicate-flag helper in B5 from text ("dup"/"first") to numeric
# (0/1) so it composes with SUMIF. Column B/D are also widened now that
# "vacant" and "unavailable" get their own columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: duplicate-unit-type flag now numeric (0/1) instead of text ---
$ws.Range("B5").Value = '&=&=IF(A5=A{-1}, 0,1)'

# --- Row 6: totals, realigned one column to the right + D's SUMIF source ---
$ws.Range("D6").Value = '&=&=SUMIF(B5:B{-1}, 1, D5:D{-1})'
$ws.Range("E6").Value = '&=&=SUM(E5:E{-1})'
$ws.Range("F6").Value = '&=&=SUM(F5:F{-1})'
$ws.Range("G6").Value = '&=&=SUM(G5:G{-1})'
$ws.Range("H6").Value = '&=&=SUM(H5:H{-1})'
$ws.Range("I6").Value = '&=&=SUM(I5:I{-1})'
$ws.Range("J6").Value = '&=&=SUM(J5:J{-1})'
$ws.Range("K6").Value = '&=&=SUM(K5:K{-1})'
$ws.Range("L6").Value = '&=&=SUM(L5:L{-1})'
$ws.Range("M6").Value = '&=&=SUM(M5:M{-1})'
$ws.Range("N6").Value = '&=&=SUM(N5:N{-1})/100'

# --- Column widths: split A:B into distinct widths, widen B & D ---
# (column A keeps its original width - only touch B and D so A gets its
# own <col> entry instead of staying merged with B's span)
$ws.Columns.Item(2).ColumnWidth = 22.5
$ws.Columns.Item(4).ColumnWidth = 27.8

# --- Selection moves to O11 ---
[void]$ws.Range("O11").Select()
